$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.766.29'
$ws.Range("E2").Value = '  +0.61%  '
$ws.Range("D3").Value = '2.508.22'
$ws.Range("E3").Value = '  +0.19%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = "'576.23"
$ws.Range("E5").Value = '  +0.15%  '
$ws.Range("D6").Value = "'167.42"
$ws.Range("E6").Value = '  +0.87%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").Value = '2.507.62'
$ws.Range("E9").Value = '  +0.16%  '
$ws.Range("D10").Value = "'0.162"
$ws.Range("E11").Value = '  +0.09%  '
$ws.Range("E12").Value = '  +4.35%  '
$ws.Range("D13").Value = "'4.95"
$ws.Range("E13").Value = '  +2.38%  '
$ws.Range("D14").Value = '2.967.69'
$ws.Range("E14").Value = '  +0.16%  '
$ws.Range("E15").Value = '  +3.22%  '
$ws.Range("D16").Value = '69.564.73'
$ws.Range("E16").Value = '  +0.54%  '
$ws.Range("E17").Value = '  +0.98%  '
$ws.Range("D18").Value = '2.489.26'
$ws.Range("E18").Value = '  -0.62%  '
$ws.Range("D19").Value = "'11.25"
$ws.Range("E19").Value = '  -0.92%  '
$ws.Range("D20").Value = "'7.50"
$ws.Range("E20").Value = '  -3.20%  '
$ws.Range("D21").Value = "'349.59"
$ws.Range("E21").Value = '  +0.98%  '
$ws.Range("E22").Value = '  -0.09%  '
$ws.Range("E23").Value = '  +1.04%  '
$ws.Range("E24").Value = '  -0.22%  '
$ws.Range("D25").Value = "'70.55"
$ws.Range("E25").Value = '  +3.45%  '
$ws.Range("D26").Value = "'3.96"
$ws.Range("E26").Value = '  +0.40%  '
$ws.Range("D27").Value = "'8.82"
$ws.Range("E27").Value = '  -0.44%  '
$ws.Range("D28").Value = '2.592.22'
$ws.Range("E28").Value = '  -1.57%  '
$ws.Range("D29").Value = "'0.996"
$ws.Range("E29").Value = '  -0.90%  '
$ws.Range("D30").Value = '0.0₃0894'
$ws.Range("E30").Value = '  +0.31%  '
$ws.Range("D31").Value = "'7.85"
$ws.Range("E31").Value = '  +0.17%  '
$ws.Range("D32").Value = "'460.56"
$ws.Range("E32").Value = '  -0.29%  '
$ws.Range("E33").Value = '  -2.45%  '
$ws.Range("E34").Value = '  -0.29%  '
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("E36").Value = '  +2.08%  '
$ws.Range("D37").Value = "'158.14"
$ws.Range("E37").Value = '  +3.41%  '
$ws.Range("D39").Value = "'18.50"
$ws.Range("E41").Value = '  +1.62%  '
$ws.Range("E42").Value = '  -0.21%  '
$ws.Range("E43").Value = '  +1.01%  '
$ws.Range("D44").Value = "'38.08"
$ws.Range("E44").Value = '  +0.13%  '
$ws.Range("E45").Value = '  -3.23%  '
$ws.Range("E46").Value = '  -6.33%  '
$ws.Range("D47").Value = "'141.57"
$ws.Range("E47").Value = '  -0.48%  '
$ws.Range("E48").Value = '  -0.11%  '
$ws.Range("E49").Value = '  -0.66%  '
$ws.Range("E50").Value = '  +0.91%  '
$ws.Range("D51").Value = "'0.580"
$ws.Range("E51").Value = '  -0.59%  '
